# Scheduled-runner style refresh of market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets. Leve metadata columns (A-G) are
# untouched; only the market-data-derived numeric cells (H-N) are updated.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 3564.3
$ws.Cells.Item(86, 9).Value = 3140.5
$ws.Cells.Item(86, 11).Value = 3140.5
$ws.Cells.Item(86, 13).Value = -2017.5

$ws.Cells.Item(89, 8).Value = 3564.3
$ws.Cells.Item(89, 9).Value = 3140.5
$ws.Cells.Item(89, 11).Value = 15702.5
$ws.Cells.Item(89, 13).Value = -10086.5

$ws.Cells.Item(107, 8).Value = 3972.303
$ws.Cells.Item(107, 9).Value = 3384.3928
$ws.Cells.Item(107, 11).Value = 3384.3928
$ws.Cells.Item(107, 13).Value = -1464.3928

$ws.Cells.Item(112, 8).Value = 2694.9666
$ws.Cells.Item(112, 10).Value = 3025.7
$ws.Cells.Item(112, 12).Value = 9077.099999999999
$ws.Cells.Item(112, 14).Value = -11293.1

$ws.Cells.Item(131, 8).Value = 12733.846
$ws.Cells.Item(131, 9).Value = 1295
$ws.Cells.Item(131, 11).Value = 3885
$ws.Cells.Item(131, 13).Value = 1155

$ws.Cells.Item(132, 8).Value = 9546.393
$ws.Cells.Item(132, 9).Value = 10470.833
$ws.Cells.Item(132, 10).Value = 3999.75
$ws.Cells.Item(132, 11).Value = 31412.499
$ws.Cells.Item(132, 12).Value = 11999.25
$ws.Cells.Item(132, 13).Value = -28882.499
$ws.Cells.Item(132, 14).Value = -17059.25

$ws.Cells.Item(133, 8).Value = 84999
$ws.Cells.Item(133, 10).Value = 84999
$ws.Cells.Item(133, 12).Value = 84999
$ws.Cells.Item(133, 14).Value = -95119

$ws.Cells.Item(134, 8).Value = 86999.5
$ws.Cells.Item(134, 10).Value = 86999.5
$ws.Cells.Item(134, 12).Value = 86999.5
$ws.Cells.Item(134, 14).Value = -97139.5

$ws.Cells.Item(139, 8).Value = 76939.60000000001
$ws.Cells.Item(139, 10).Value = 76939.60000000001
$ws.Cells.Item(139, 12).Value = 76939.60000000001
$ws.Cells.Item(139, 14).Value = -87219.60000000001

$ws.Cells.Item(140, 8).Value = 107933
$ws.Cells.Item(140, 10).Value = 107933
$ws.Cells.Item(140, 12).Value = 107933
$ws.Cells.Item(140, 14).Value = -118293

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1820.4445
$ws.Cells.Item(45, 9).Value = 1730.6666
$ws.Cells.Item(45, 10).Value = 2000
$ws.Cells.Item(45, 11).Value = 1730.6666
$ws.Cells.Item(45, 12).Value = 2000
$ws.Cells.Item(45, 13).Value = -1353.6666
$ws.Cells.Item(45, 14).Value = -2754

$ws.Cells.Item(74, 8).Value = 2530720.8
$ws.Cells.Item(74, 9).Value = 5053533
$ws.Cells.Item(74, 11).Value = 5053533
$ws.Cells.Item(74, 13).Value = -5052659

$ws.Cells.Item(77, 8).Value = 2530720.8
$ws.Cells.Item(77, 9).Value = 5053533
$ws.Cells.Item(77, 11).Value = 25267665
$ws.Cells.Item(77, 13).Value = -25263297

$ws.Cells.Item(110, 8).Value = 1612.8334
$ws.Cells.Item(110, 9).Value = 1372
$ws.Cells.Item(110, 11).Value = 1372
$ws.Cells.Item(110, 13).Value = 673

$ws.Cells.Item(122, 8).Value = 1411.0769
$ws.Cells.Item(122, 9).Value = 1262.909
$ws.Cells.Item(122, 11).Value = 3788.727
$ws.Cells.Item(122, 13).Value = -1338.727

$ws.Cells.Item(134, 8).Value = 50475.2
$ws.Cells.Item(134, 10).Value = 50475.2
$ws.Cells.Item(134, 12).Value = 50475.2
$ws.Cells.Item(134, 14).Value = -60615.2

$ws.Cells.Item(139, 8).Value = 120715
$ws.Cells.Item(139, 10).Value = 120715
$ws.Cells.Item(139, 12).Value = 120715
$ws.Cells.Item(139, 14).Value = -130995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 42500
$ws.Cells.Item(26, 9).Value = 20000
$ws.Cells.Item(26, 11).Value = 20000
$ws.Cells.Item(26, 13).Value = -19708

$ws.Cells.Item(96, 8).Value = 41182.668
$ws.Cells.Item(96, 9).Value = 31786.5
$ws.Cells.Item(96, 11).Value = 31786.5
$ws.Cells.Item(96, 13).Value = -29040.5

$ws.Cells.Item(105, 8).Value = 2472.1191
$ws.Cells.Item(105, 9).Value = 1984.8387
$ws.Cells.Item(105, 11).Value = 1984.8387
$ws.Cells.Item(105, 13).Value = -237.8387

$ws.Cells.Item(134, 8).Value = 4766777.5
$ws.Cells.Item(134, 9).Value = 4694.839
$ws.Cells.Item(134, 11).Value = 14084.517
$ws.Cells.Item(134, 13).Value = -11549.517

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 41384.5
$ws.Cells.Item(99, 9).Value = 41383.168
$ws.Cells.Item(99, 10).Value = 41388.5
$ws.Cells.Item(99, 11).Value = 41383.168
$ws.Cells.Item(99, 12).Value = 41388.5
$ws.Cells.Item(99, 13).Value = -39885.168
$ws.Cells.Item(99, 14).Value = -44384.5

$ws.Cells.Item(126, 8).Value = 41384.5
$ws.Cells.Item(126, 9).Value = 41383.168
$ws.Cells.Item(126, 10).Value = 41388.5
$ws.Cells.Item(126, 11).Value = 124149.504
$ws.Cells.Item(126, 12).Value = 124165.5
$ws.Cells.Item(126, 13).Value = -121679.504
$ws.Cells.Item(126, 14).Value = -129105.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 17960906
$ws.Cells.Item(4, 9).Value = 26939924
$ws.Cells.Item(4, 11).Value = 80819772
$ws.Cells.Item(4, 13).Value = -80819660

$ws.Cells.Item(68, 8).Value = 5162.143
$ws.Cells.Item(68, 10).Value = 5162.143
$ws.Cells.Item(68, 12).Value = 15486.429
$ws.Cells.Item(68, 14).Value = -17108.429

$ws.Cells.Item(71, 8).Value = 5162.143
$ws.Cells.Item(71, 10).Value = 5162.143
$ws.Cells.Item(71, 12).Value = 46459.287
$ws.Cells.Item(71, 14).Value = -54571.287

$ws.Cells.Item(87, 8).Value = 22000
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 13).ClearContents()

$ws.Cells.Item(90, 8).Value = 22000
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(90, 13).ClearContents()

$ws.Cells.Item(123, 8).Value = 13171.667
$ws.Cells.Item(123, 9).Value = 9515
$ws.Cells.Item(123, 11).Value = 28545
$ws.Cells.Item(123, 13).Value = -26095

$ws.Cells.Item(126, 8).Value = 8338.111000000001
$ws.Cells.Item(126, 9).Value = 2208.6
$ws.Cells.Item(126, 11).Value = 6625.799999999999
$ws.Cells.Item(126, 13).Value = -1685.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 41513.633
$ws.Cells.Item(70, 9).Value = 10329.294
$ws.Cells.Item(70, 11).Value = 10329.294
$ws.Cells.Item(70, 13).Value = -10059.294

$ws.Cells.Item(73, 8).Value = 41513.633
$ws.Cells.Item(73, 9).Value = 10329.294
$ws.Cells.Item(73, 11).Value = 10329.294
$ws.Cells.Item(73, 13).Value = -9393.294

$ws.Cells.Item(80, 8).Value = 5821.9375
$ws.Cells.Item(80, 9).Value = 3633.0833
$ws.Cells.Item(80, 11).Value = 3633.0833
$ws.Cells.Item(80, 13).Value = -2635.0833

$ws.Cells.Item(83, 8).Value = 5821.9375
$ws.Cells.Item(83, 9).Value = 3633.0833
$ws.Cells.Item(83, 11).Value = 18165.4165
$ws.Cells.Item(83, 13).Value = -13173.4165

$ws.Cells.Item(102, 8).Value = 2358.5
$ws.Cells.Item(102, 9).Value = 2287.2222
$ws.Cells.Item(102, 11).Value = 2287.2222
$ws.Cells.Item(102, 13).Value = -665.2222000000002

$ws.Cells.Item(132, 8).Value = 4288.9395
$ws.Cells.Item(132, 9).Value = 3817.2808
$ws.Cells.Item(132, 10).Value = 7276.1113
$ws.Cells.Item(132, 11).Value = 11451.8424
$ws.Cells.Item(132, 12).Value = 21828.3339
$ws.Cells.Item(132, 13).Value = -8921.8424
$ws.Cells.Item(132, 14).Value = -26888.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4699
$ws.Cells.Item(40, 9).Value = 4854.778
$ws.Cells.Item(40, 11).Value = 4854.778
$ws.Cells.Item(40, 13).Value = -4718.778

$ws.Cells.Item(55, 8).Value = 1617.15
$ws.Cells.Item(55, 9).Value = 1985.2858
$ws.Cells.Item(55, 10).Value = 1418.9231
$ws.Cells.Item(55, 11).Value = 1985.2858
$ws.Cells.Item(55, 12).Value = 1418.9231
$ws.Cells.Item(55, 13).Value = -1812.2858
$ws.Cells.Item(55, 14).Value = -1764.9231

$ws.Cells.Item(68, 8).Value = 1722.45
$ws.Cells.Item(68, 10).Value = 1400
$ws.Cells.Item(68, 12).Value = 1400
$ws.Cells.Item(68, 14).Value = -2898

$ws.Cells.Item(71, 8).Value = 1722.45
$ws.Cells.Item(71, 10).Value = 1400
$ws.Cells.Item(71, 12).Value = 7000
$ws.Cells.Item(71, 14).Value = -14488

$ws.Cells.Item(132, 8).Value = 6153079
$ws.Cells.Item(132, 9).Value = 11687388
$ws.Cells.Item(132, 11).Value = 35062164
$ws.Cells.Item(132, 13).Value = -35059634

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4067550.5
$ws.Cells.Item(132, 9).Value = 4631934
$ws.Cells.Item(132, 10).Value = 3989.8
$ws.Cells.Item(132, 11).Value = 13895802
$ws.Cells.Item(132, 12).Value = 11969.4
$ws.Cells.Item(132, 13).Value = -13893272
$ws.Cells.Item(132, 14).Value = -17029.4
